# ---------------------------------------------------------------------------
# Translate the two existing algorithm headers to Vietnamese, add a third
# "Di truyền" (Genetic) algorithm block (columns L:P) mirroring the existing
# "Tham lam" / "Quy hoạch động" blocks, and extend the results table with
# eight more problem rows (3-10), filling in the new L:P columns for the
# rows that already existed (Problem 1 & 2) as well.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Translate the two existing merged header labels -----------------------
$ws.Range("B1").Value = "Tham lam"
$ws.Range("G1").Value = "Quy hoạch động"

# --- New header block for the "Di truyền" algorithm (L1:P1) ---------------
# Merge first, then paste the formatting from the existing "Tham lam" header
# so the new block picks up the same bold/border/center style (s=1).
$ws.Range("L1:P1").Merge()
$ws.Range("B1").Copy()
$ws.Range("L1:P1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("L1").Value = "Di truyền"

# --- New row-2 statistic sub-headers for L:P, copied from G2:K2 -----------
$ws.Range("G2:K2").Copy()
$ws.Range("L2:P2").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$statNames = @("mean", "std", "min", "med", "max")
$statCols  = @("L", "M", "N", "O", "P")
for ($i = 0; $i -lt $statNames.Length; $i++) {
    $ws.Range(($statCols[$i] + "2")).Value = $statNames[$i]
}

# --- Data rows (columns B:P) for every problem, keyed by row number -------
$dataCols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P")

$rows = @{
    4  = @{ Name = "Problem 1";  Values = @(2758, 184, 2554, 2690, 3047, 181, 1, 179, 181, 183, 11932, 894, 11213, 11587, 14025) }
    5  = @{ Name = "Problem 2";  Values = @(2810, 263, 2641, 2689, 3467, 109, 2, 107, 109, 114, 11677, 744, 10715, 11404, 12970) }
    6  = @{ Name = "Problem 3";  Values = @(1438, 1383, 13, 1428, 2945, 157, 1, 156, 157, 158, 5489, 303, 5252, 5397, 6351) }
    7  = @{ Name = "Problem 4";  Values = @(2, 1, 1, 2, 3, 91, 1, 89, 90, 94, 5458, 86, 5310, 5460, 5652) }
    8  = @{ Name = "Problem 5";  Values = @(2, 2, 1, 2, 6, 78, 2, 76, 77, 83, 3469, 265, 3237, 3394, 4224) }
    9  = @{ Name = "Problem 6";  Values = @(2, 1, 1, 2, 3, 59, 1, 57, 59, 63, 4232, 177, 4054, 4209, 4697) }
    10 = @{ Name = "Problem 7";  Values = @(1, 0, 1, 1, 2, 115, 1, 113, 115, 116, 3156, 48, 3096, 3147, 3278) }
    11 = @{ Name = "Problem 8";  Values = @(2, 0, 1, 2, 3, 49, 0, 49, 49, 50, 2813, 76, 2701, 2782, 2930) }
    12 = @{ Name = "Problem 9";  Values = @(5, 1, 3, 5, 7, 346, 2, 344, 345, 352, 8527, 256, 8242, 8429, 9163) }
    13 = @{ Name = "Problem 10"; Values = @(2, 1, 1, 2, 3, 54, 1, 53, 54, 57, 4026, 66, 3926, 4017, 4160) }
}

# Existing rows (4 & 5) only need their B:P values updated/extended.
# Brand-new rows (6-13) also need the "Problem N" label in column A, styled
# like the existing problem-label cells (copy format from row 5's A cell).
foreach ($r in (4..13)) {
    $row = $rows[$r]

    if ($r -ge 6) {
        $ws.Range("A5").Copy()
        $ws.Range(("A" + $r)).PasteSpecial(-4122)   # xlPasteFormats
        $excel.CutCopyMode = $false
        $ws.Range(("A" + $r)).Value = $row.Name
    }

    for ($i = 0; $i -lt $dataCols.Length; $i++) {
        $ws.Range(($dataCols[$i] + $r)).Value = $row.Values[$i]
    }
}
